$d = $word.ActiveDocument

# --- 1) Move the "_GoBack" last-edit bookmark ---------------------------
# It currently sits between "2008-03 " and "Document " in the title.
# Remove it there, and re-insert it inside the run that is split by the
# edit "which produced version  2008-03-11" -> "which pro" | "duced version  2008-03-11".
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$found = $d.Content
$found.Find.Execute("which produced version  2008-03-11", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$splitPoint = $found.Start + [string]"which pro".Length

$gobackRange = $d.Range($splitPoint, $splitPoint)
$d.Bookmarks.Add("_GoBack", $gobackRange)

# --- 2) Set the Normal style font to Calibri 11pt (non-heading text) ----
$normal = $d.Styles("Normal")
$normal.Font.Name = "Calibri"
$normal.Font.Size = 11
